# feat: update foreach/endrow/endloop with new behaviour
#
# The "#! END_ROW" template directive gains an optional boolean argument.
# A new column K is added next to the existing "#! END_ROW" marker column
# (J) to show/exercise the new directive variant:
#   K2 -> "#! END_ROW true"   (new directive variant -> new shared string)
#   K3 -> "#! END_ROW"        (plain directive, reuses the existing string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "#! END_ROW true"
$ws.Range("K3").Value = "#! END_ROW"

# Leave the view scrolled back to the left edge with the new K3 cell
# selected, matching the template author's final view state.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("K3").Select() | Out-Null
